$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B9").Select()
